$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new rightmost data column (S) for year 2023, matching
# the formatting already used by the neighbouring 2022 column (R).
# Copy R4:R6's formatting into S4:S6 first, then overwrite with the new values
# so number formats / fonts / borders / alignment for the new column match
# the existing table exactly.
$ws.Range("R4:R6").Copy()
$ws.Range("S4:S6").PasteSpecial(-4122)

$ws.Range("S4").Value = 2023
$ws.Range("S5").Value = 7.1262361838278068
$ws.Range("S6").Value = 10.974456007568591
